# Initial missing file exception fix
# Adds BURHILL SPORTSW / AGDISC / CROWN price columns (S, T, U) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cells (new shared strings: SPORTSW, AGDISC, CROWN)
$ws.Range("S1").Value = "SPORTSW"
$ws.Range("T1").Value = "AGDISC"
$ws.Range("U1").Value = "CROWN"

# New data columns S (=BURHILL, rounded), T (discounted), U (=BURHILL, rounded)
$ws.Range("S2").Value = 130.5
$ws.Range("T2").Value = 126.15
$ws.Range("U2").Value = 130.5
$ws.Range("S3").Value = 92.25
$ws.Range("T3").Value = 89.18
$ws.Range("U3").Value = 92.25
$ws.Range("S4").Value = 95.4
$ws.Range("T4").Value = 92.22
$ws.Range("U4").Value = 95.4
$ws.Range("S5").Value = 86.4
$ws.Range("T5").Value = 83.52
$ws.Range("U5").Value = 86.4
$ws.Range("S6").Value = 76.5
$ws.Range("T6").Value = 73.95
$ws.Range("U6").Value = 76.5
$ws.Range("S7").Value = 64.349999999999994
$ws.Range("T7").Value = 62.21
$ws.Range("U7").Value = 64.349999999999994
$ws.Range("S8").Value = 61.92
$ws.Range("T8").Value = 59.86
$ws.Range("U8").Value = 61.92
$ws.Range("S9").Value = 42.75
$ws.Range("T9").Value = 41.33
$ws.Range("U9").Value = 42.75
$ws.Range("S10").Value = 0
$ws.Range("T10").Value = 0
$ws.Range("U10").Value = 0
$ws.Range("S11").Value = 0
$ws.Range("T11").Value = 0
$ws.Range("U11").Value = 0
$ws.Range("S12").Value = 53.55
$ws.Range("T12").Value = 51.77
$ws.Range("U12").Value = 53.55
$ws.Range("S13").Value = 54
$ws.Range("T13").Value = 52.2
$ws.Range("U13").Value = 54
$ws.Range("S14").Value = 47.25
$ws.Range("T14").Value = 45.68
$ws.Range("U14").Value = 47.25
$ws.Range("S15").Value = 41.4
$ws.Range("T15").Value = 40.02
$ws.Range("U15").Value = 41.4
$ws.Range("S16").Value = 25.65
$ws.Range("T16").Value = 24.8
$ws.Range("U16").Value = 25.65
$ws.Range("S17").Value = 18.899999999999999
$ws.Range("T17").Value = 18.27
$ws.Range("U17").Value = 18.899999999999999
$ws.Range("S18").Value = 33.75
$ws.Range("T18").Value = 32.630000000000003
$ws.Range("U18").Value = 33.75
$ws.Range("S19").Value = 30.15
$ws.Range("T19").Value = 29.15
$ws.Range("U19").Value = 30.15
$ws.Range("S20").Value = 23.4
$ws.Range("T20").Value = 22.62
$ws.Range("U20").Value = 23.4
$ws.Range("S21").Value = 23.4
$ws.Range("T21").Value = 22.62
$ws.Range("U21").Value = 23.4
$ws.Range("S22").Value = 20.7
$ws.Range("T22").Value = 20.010000000000002
$ws.Range("U22").Value = 20.7
$ws.Range("S23").Value = 20.7
$ws.Range("T23").Value = 20.010000000000002
$ws.Range("U23").Value = 20.7
$ws.Range("S24").Value = 27
$ws.Range("T24").Value = 26.1
$ws.Range("U24").Value = 27
$ws.Range("S25").Value = 64.349999999999994
$ws.Range("T25").Value = 62.21
$ws.Range("U25").Value = 64.349999999999994
$ws.Range("S26").Value = 42.75
$ws.Range("T26").Value = 41.33
$ws.Range("U26").Value = 42.75
$ws.Range("S27").Value = 25.65
$ws.Range("T27").Value = 24.8
$ws.Range("U27").Value = 25.65
$ws.Range("S28").Value = 18.899999999999999
$ws.Range("T28").Value = 18.27
$ws.Range("U28").Value = 18.899999999999999
$ws.Range("S29").Value = 45
$ws.Range("T29").Value = 43.5
$ws.Range("U29").Value = 45
$ws.Range("S30").Value = 36
$ws.Range("T30").Value = 34.799999999999997
$ws.Range("U30").Value = 36

# Restore the final selection left by the editor
$ws.Range("S18").Select()
